# Weekly update: insert a new Berenjena price record at the top of the
# "Vega Modelo de Temuco" data block (row 204), shifting the existing
# rows 204-218 down to 205-219.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 204, pushing the rows below it down
# (mirrors Excel's Rows("204:204").Insert Shift:=xlShiftDown).
$ws.Rows.Item(204).Insert(-4121)

# Populate the freshly inserted row with the new weekly record.
$ws.Cells.Item(204, 1).Value = 10
$ws.Cells.Item(204, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(204, 3).Value = 'La Araucanía'
$ws.Cells.Item(204, 4).Value = 44578
$ws.Cells.Item(204, 5).Value = 9
$ws.Cells.Item(204, 6).Value = 100112001
$ws.Cells.Item(204, 7).Value = 'Berenjena'
$ws.Cells.Item(204, 8).Value = 'Sin especificar'
$ws.Cells.Item(204, 9).Value = 'Segunda'
$ws.Cells.Item(204, 10).Value = 30
$ws.Cells.Item(204, 11).Value = 10000
$ws.Cells.Item(204, 12).Value = 10000
$ws.Cells.Item(204, 13).Value = 10000
$ws.Cells.Item(204, 14).Value = '$/caja 90 unidades'
$ws.Cells.Item(204, 15).Value = 'Región del Maule'
$ws.Cells.Item(204, 16).Value = 111
$ws.Cells.Item(204, 17).Value = 90
$ws.Cells.Item(204, 18).Value = 'Hortaliza'
